# update database and change read_price algorithm
# Shift each yearly income-statement column one period to the left
# (drop the oldest period, 1396/12, and append the newest, 1401/12),
# plus a handful of corrected figures (the "read_price algorithm" fix).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Period / publish-date headers (row 8 / row 9, columns D:H) ----
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

$ws.Range("D9").Value = "1399-04-19 (8)"
$ws.Range("E9").Value = "1400-04-26 (10)"
$ws.Range("F9").Value = "1401-04-21 (10)"
$ws.Range("G9").Value = "1402-02-29 (8)"
$ws.Range("H9").Value = "1402-02-29"

# ---- Income statement data rows (row 11 .. row 27), columns D:H ----
$data = @{
    11 = @(485390122, 866331631, 1109548371, 2130779345, 3397308491)
    12 = @(-429882294, -811599243, -981308661, -1936806043, -2952758372)
    13 = @(55507828, 54732388, 128239710, 193973302, 444550119)
    14 = @(-5212276, -7262785, -8311647, -12618898, -26169707)
    15 = @(0, 0, 0, 0, 0)
    16 = @(9400897, -19881449, 38672654, 11984547, -221693)
    17 = @(59696449, 27588154, 158600717, 193338951, 418158719)
    18 = @(0, 0, -88767, -6545203, -16411649)
    19 = @(164963, 4136439, 7280248, 8061224, 25883616)
    20 = @(59861412, 31724593, 165792198, 194854972, 427630686)
    21 = @(-9086856, -3199962, -17342856, -29274798, -45768701)
    22 = @(50774556, 28524631, 148449342, 165580174, 381861985)
    23 = @(0, 0, 0, 0, 0)
    24 = @(50774556, 28524631, 148449342, 165580174, 381861985)
    25 = @(996, 0, 1071, 871, 1682)
    26 = @(51000000, 51000000, 138649318, 190000000, 227000000)
    27 = @(224, 126, 654, 729, 1682)
}

$cols = @("D", "E", "F", "G", "H")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value = $vals[$i]
    }
}
